# Product Requisition workbook update — "15-1024 report"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits on Sheet1 ---
# Row 9  (MMST 470): clear quantity
$ws.Range("C9").Value = ""

# Row 14 (Swap SIM): clear quantity
$ws.Range("C14").Value = ""

# Row 31 (SC Voice-19): quantity 15000 -> 10000
$ws.Range("C31").Value = 10000

# Row 32 (SC Data -29): clear quantity
$ws.Range("C32").Value = ""

# Row 33 (SC Voice -29): quantity blank -> 1000
$ws.Range("C33").Value = 1000

# Row 43 (I'top-Up): quantity 472416 -> 158961
$ws.Range("C43").Value = 158961

# --- View-state: scroll position + active selection ---
$ws.Activate()
$ws.Range("H51").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
